# Applies the "Celestial Symphony" -> "A Journey Through Biology" rewrite
# (title/author/email/body/summary text swap + font-name fix + a trailing
# empty paragraph) to $word.ActiveDocument.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$old,
        [string]$new
    )
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.Text = $old
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = $new
    $null = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---------------------------------------------------------------------------
# 1) Fix the (misspelled) font name everywhere. Using a Range that excludes
#    the very last (paragraph-mark) position keeps Word from stamping a new
#    paragraph-mark rPr on every paragraph.
# ---------------------------------------------------------------------------
$fontRange = $d.Range(0, $d.Content.End - 1)
$fontRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# 2) Title / author / email swap.
# ---------------------------------------------------------------------------
Replace-Text "Celestial Symphony: Unveiling the Harmony of the Cosmos" "A Journey Through Biology: Unraveling the Wonders of Life"
Replace-Text "Evelyn Abernathy" "Valerie Knight"
Replace-Text "evelyn" "valerieknight@emailcentral"
Replace-Text "abernathy@stardustobservatory." ""
Replace-Text "org" "net"

# ---------------------------------------------------------------------------
# 3) Body paragraph (the long "celestial symphony" essay -> biology essay).
# ---------------------------------------------------------------------------
Replace-Text "In the vast expanse of the universe, there exists a celestial symphony--a harmonious blend of rhythmic cosmic phenomena" "Biology, the study of life, invites us on an enthralling journey to unravel the wonders of the living world"
Replace-Text " Unraveling this symphony requires blending disciplines like astronomy, physics, and mathematics, revealing the intricate dance of cosmic bodies and the underlying forces that govern the universe" " From the vastness of an African savanna to the intricate inner workings of a human cell, biology unveils a tapestry of interconnectedness and diversity that inspires awe and curiosity"
Replace-Text " From the rhythmic pulsations of celestial bodies to the ethereal ebb and flow of gravitational waves, the universe resonates with an intriguing melody, inviting us to decode its secrets" " As we delve into this fascinating discipline, we embark on a quest to comprehend the fundamental principles governing life, from the genetic code that orchestrates biological processes to the evolutionary forces shaping the history of organisms"
Replace-Text "The cosmic dance of stellar evolution, from the birth of stars to their final metamorphosis into remnants, echoes the rhythm of life and decay" "Our exploration begins by peering into the microscopic realm, where atoms and molecules dance in intricate patterns to form the building blocks of life"
Replace-Text " The harmonious interplay of galaxies, clusters, and superclusters, bound together by the invisible threads of gravity, mirrors the interconnectedness of life on Earth" " We marvel at the complexity of cells, the fundamental units of living organisms, and discover the specialized structures and functions that enable them to thrive"
Replace-Text " Within this celestial symphony, the pulsating brilliance of pulsars and the graceful ballet of orbiting celestial bodies add their unique melodies, underscoring the complexity and unity of the universe" " We investigate the processes of cell division, metabolism, and genetics, understanding how cells grow, reproduce, and pass on their traits to future generations"
Replace-Text "The harmony of the cosmos extends beyond the visible sphere" "Moving beyond the cellular level, we delve into the intricate world of organisms, examining their adaptations, behaviors, and interactions within ecosystems"
Replace-Text " The silent whispers of gravitational waves carry tales of cosmic events, echoing the cataclysmic dance of colliding black holes or the exhilarating merger of neutron stars" " We learn about the incredible diversity of life on Earth, from the vibrant colors of coral reefs to the soaring heights of mountain gorillas"
Replace-Text " These gravitational waves, like cosmic drumbeats, paint a vivid picture of the universe's energetic tapestry, revealing hidden dimensions of reality" " We explore the delicate balance of ecosystems, appreciating the intricate relationships between species and their environment, and recognize the crucial role humans play in preserving and protecting the natural world"
Replace-Text ". As we unravel the intricate notes of this celestial symphony, we gain deeper insights into the fundamental laws that orchestrate the universe" ""

# ---------------------------------------------------------------------------
# 4) Summary paragraph.
# ---------------------------------------------------------------------------
Replace-Text "The universe is a resonant symphony of cosmic phenomena, interwoven with the rhythmic dances of celestial bodies and the ethereal melodies of gravitational waves" "Through our journey into biology, we have gained a deep understanding of life's intricacies, from the molecular foundations of cells to the vast interconnectedness of ecosystems"
Replace-Text " This interplay unveils intricate patterns of harmony, connecting the vast expanse of the universe through its inherent interconnectedness" " We have explored the fascinating processes that govern biological systems, unlocking the secrets of genetics, evolution, and adaptation"
Replace-Text " As we decipher this celestial symphony, we enrich our understanding of the fundamental laws governing the universe, expanding our knowledge of its origins and evolution" " Above all, we have developed a profound appreciation for the beauty and complexity of life in all its forms, inspiring us to protect and preserve the natural world for generations to come"
Replace-Text ". Through this cosmic exploration, we find echoes of unity and harmony, blurring the boundaries between disciplines and transcending earthly limitations" ""

# ---------------------------------------------------------------------------
# 5) A new trailing empty paragraph appears right before the section break.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$endOfBody = $lastPara.Range.End
$insertion = $d.Range($endOfBody - 1, $endOfBody - 1)
$insertion.InsertParagraphAfter()
